# Actualización automática 2025-10-06 15:30:20
# Target sheet: "CUMPLIMIENTO MENSUAL" (3rd worksheet) of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Column width adjustments (D, E, F) ---
# Excel's ColumnWidth (character units) is offset from the stored OOXML
# <col width="..."> by a constant padding of 5/6 of a character for the
# workbook's default font, so subtract that to land exactly on the target
# stored widths of 11 / 22 / 18.
$ws.Columns.Item(4).ColumnWidth = 11 - (5/6)
$ws.Columns.Item(5).ColumnWidth = 22 - (5/6)
$ws.Columns.Item(6).ColumnWidth = 18 - (5/6)

# --- Row 3: 240X80 PORCELANATO -> VENTA reset to 0 ---
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 7475.4083879616
$ws.Range("F3").Value = 0

# --- Row 4: FREGADEROS DE COCINA -> VENTA reset to 0 ---
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 485.098641648355
$ws.Range("F4").Value = 0

# --- Row 5: GRIFERIAS -> VENTA reset to 0 ---
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 150
$ws.Range("F5").Value = 0

# --- Row 6: INODOROS -> VENTA reset to 0 ---
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 907.166108615601
$ws.Range("F6").Value = 0

# --- Row 9: OTROS -> VENTA reset to 0 ---
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0

# --- Row 11: PIEDRA SINTERIZADA -> VENTA reset to 0 ---
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 2922.22458185274
$ws.Range("F11").Value = 0

# --- Row 12: PORCELANATO -> PRESUPUESTO updated, VENTA reset to 0 ---
$ws.Range("C12").Value = 27954.98
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 27954.98
$ws.Range("F12").Value = 0

# --- Remove row 14 (SAL SOLUBLE) entirely; TOTAL row (old row 15) shifts up to row 14 ---
$ws.Rows.Item(14).Delete()

# --- Update the new row 14 (TOTAL) to reflect the new column sums ---
$ws.Range("C14").Value = 42203.38110009469
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 42203.38110009469
$ws.Range("F14").Value = 0
